$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D2:D51 to text format so numeric-looking strings (e.g. "1.000", "125.10")
# are preserved exactly as text rather than being normalized into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.824.02'
$ws.Range("E2").Value = '  -1.04%  '
$ws.Range("D3").Value = '1.763.40'
$ws.Range("E3").Value = '  -2.61%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '321.54'
$ws.Range("E5").Value = '  -2.56%  '
$ws.Range("D6").Value = '0.9986'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("D7").Value = '0.4249'
$ws.Range("E7").Value = '  -4.35%  '
$ws.Range("D8").Value = '0.3634'
$ws.Range("E8").Value = '  -2.46%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D9").Value = '42.53'
$ws.Range("E9").Value = '  -4.86%  '
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.07503'
$ws.Range("E10").Value = '  -2.60%  '
$ws.Range("D11").Value = '1.094'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  +0.01%  '
$ws.Range("D13").Value = '20.71'
$ws.Range("E13").Value = '  -5.73%  '
$ws.Range("D14").Value = '6.070'
$ws.Range("E14").Value = '  -3.67%  '
$ws.Range("D15").Value = '7.282'
$ws.Range("E15").Value = '  -2.56%  '
$ws.Range("D16").Value = '1.786.23'
$ws.Range("E16").Value = '  -1.56%  '
$ws.Range("D17").Value = '91.38'
$ws.Range("E17").Value = '  -2.43%  '
$ws.Range("D18").Value = '0.00001058'
$ws.Range("E18").Value = '  -2.26%  '
$ws.Range("D19").Value = '0.06375'
$ws.Range("E19").Value = '  -1.58%  '
$ws.Range("E20").Value = '  -0.23%  '
$ws.Range("D21").Value = '17.03'
$ws.Range("E21").Value = '  -2.62%  '
$ws.Range("D22").Value = '5.909'
$ws.Range("E22").Value = '  -5.57%  '
$ws.Range("D23").Value = '27.870.29'
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("D24").Value = '11.23'
$ws.Range("E24").Value = '  -4.13%  '
$ws.Range("D25").Value = '2.107'
$ws.Range("E25").Value = '  -1.37%  '
$ws.Range("D26").Value = '157.75'
$ws.Range("E26").Value = '  +1.43%  '
$ws.Range("D27").Value = '20.24'
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("D28").Value = '1.981.96'
$ws.Range("E28").Value = '  -1.78%  '
$ws.Range("D29").Value = '2.139'
$ws.Range("E29").Value = '  -8.13%  '
$ws.Range("D30").Value = '125.10'
$ws.Range("E30").Value = '  -1.66%  '
$ws.Range("D31").Value = '1.116'
$ws.Range("E31").Value = '  -7.10%  '
$ws.Range("E32").Value = '  +0.43%  '
$ws.Range("D33").Value = '5.553'
$ws.Range("E33").Value = '  -5.09%  '
$ws.Range("D34").Value = '0.08878'
$ws.Range("E34").Value = '  -3.78%  '
$ws.Range("D35").Value = '12.26'
$ws.Range("E35").Value = '  -6.32%  '
$ws.Range("D36").Value = '0.02293'
$ws.Range("E36").Value = '  -2.17%  '
$ws.Range("D37").Value = '0.2104'
$ws.Range("D38").Value = '0.06036'
$ws.Range("E38").Value = '  -2.67%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").Value = '0.6329'
$ws.Range("E39").Value = '  -3.55%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").Value = '4.961'
$ws.Range("E40").Value = '  -4.02%  '
$ws.Range("D41").Value = '1.179'
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("D42").Value = '0.9980'
$ws.Range("E42").Value = '  -0.22%  '
$ws.Range("D43").Value = '7.866'
$ws.Range("E43").Value = '  -2.66%  '
$ws.Range("D44").Value = '1.398'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("D45").Value = '13.36'
$ws.Range("E45").Value = '  -3.95%  '
$ws.Range("E46").Value = '  -3.45%  '
$ws.Range("D47").Value = '3.684'
$ws.Range("E47").Value = '  -2.20%  '
$ws.Range("D48").Value = '1.984'
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("D49").Value = '122.92'
$ws.Range("E49").Value = '  -3.04%  '
$ws.Range("D50").Value = '1.184'
$ws.Range("E50").Value = '  +3.00%  '
$ws.Range("E51").Value = '  -2.01%  '

# Restore the default (unstyled) cell style for the D column data cells so no
# extraneous formatting is introduced, matching the original workbook styling.
$ws.Range("D2:D51").Style = "Normal"
